$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.070.75"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.858.84"
$ws.Range("E3").Value = "  +3.33%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.59"
$ws.Range("E5").Value = "  +4.51%  "
$ws.Range("E6").Value = "  +2.10%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.47"
$ws.Range("E8").Value = "  +9.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.329"
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("E10").Value = "  +2.95%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.128.17"
$ws.Range("E12").Value = "  +3.28%  "
$ws.Range("E13").Value = "  +3.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.854.25"
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("E15").Value = "  +3.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.70"
$ws.Range("E16").Value = "  +3.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.058.04"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.33"
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0798"
$ws.Range("E19").Value = "  +2.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.98"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.21"
$ws.Range("E21").Value = "  +4.21%  "
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +2.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.61"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("E26").Value = "  +31.85%  "
$ws.Range("E27").Value = "  +3.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.70"
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("E29").Value = "  +3.50%  "
$ws.Range("E30").Value = "  +3.90%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  -0.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.00"
$ws.Range("E33").Value = "  +3.62%  "
$ws.Range("E34").Value = "  +14.33%  "
$ws.Range("E35").Value = "  +23.01%  "
$ws.Range("E36").Value = "  +5.18%  "
$ws.Range("E38").Value = "  +13.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "91.71"
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0202"
$ws.Range("E40").Value = "  +6.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.353.50"
$ws.Range("E41").Value = "  +3.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.88"
$ws.Range("E42").Value = "  +4.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.34"
$ws.Range("E43").Value = "  +6.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.65"
$ws.Range("E44").Value = "  +56.35%  "
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0550"
$ws.Range("E47").Value = "  +7.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.36"
$ws.Range("E48").Value = "  +3.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.041.45"
$ws.Range("E49").Value = "  +2.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0680"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("E51").Value = "  +18.24%  "
